# "add Use Item bug"
#
# The Icon column (P) for every equip row pointed at a PNG file path
# ("UI/SteampunkUI/resource/icons/img_equip.png"). The game's item-use /
# icon-loading code expects the icon path WITHOUT the file extension, so
# drop the ".png" suffix from every Icon value in the table. Once every
# cell referencing the old string is rewritten, the now-unused
# "...img_equip.png" shared string is gone from the saved workbook and the
# new extension-less string is the one left behind.
#
# Also nudge the sheet's cursor / column-P display width, matching what
# the author's Excel session left behind when they made the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite every Icon (column P) value for the data rows (2-7): strip ".png".
$ws.Range("P2:P7").Value = "UI/SteampunkUI/resource/icons/img_equip"

# Widen column P (Icon) now that bestFit no longer applies - it's a manual width.
$ws.Columns.Item(16).ColumnWidth = 47.2

# Move the selection/active cell to K23 (cosmetic - where the editor's cursor ended up).
$ws.Range("K23").Select()
